$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.882106333333333
$ws.Range("H2").Value = 11.646319
$ws.Range("I2").Value = 0.6257373677154582
$ws.Range("J2").Value = 0.6257373677154581
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.044118333333333
$ws.Range("N2").Value = 6.132354999999999
$ws.Range("O2").Value = 0.1776005292722278
$ws.Range("P2").Value = 0.1776005292722278
$ws.Range("Q2").Value = 7.935484727916109
$ws.Range("R2").Value = 71.41936255124499
$ws.Range("S2").Value = 0.111131287691676
$ws.Range("T2").Value = 0.111131287691676
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.882106333333333
$ws.Range("H3").Value = 11.646319
$ws.Range("I3").Value = 0.6257373677154582
$ws.Range("J3").Value = 0.6257373677154581
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.059280333333334
$ws.Range("N3").Value = 21.177841
$ws.Range("O3").Value = 0.6133362746356149
$ws.Range("P3").Value = 0.6133362746356149
$ws.Range("Q3").Value = 27.40487689080878
$ws.Range("R3").Value = 246.643892017279
$ws.Range("S3").Value = 0.383787426014895
$ws.Range("T3").Value = 0.383787426014895
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.882106333333333
$ws.Range("H4").Value = 11.646319
$ws.Range("I4").Value = 0.6257373677154582
$ws.Range("J4").Value = 0.6257373677154581
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.406242333333333
$ws.Range("N4").Value = 7.218726999999999
$ws.Range("O4").Value = 0.2090631960921573
$ws.Range("P4").Value = 0.2090631960921573
$ws.Range("Q4").Value = 9.34128860176811
$ws.Range("R4").Value = 84.07159741591299
$ws.Range("S4").Value = 0.1308186540088871
$ws.Range("T4").Value = 0.1308186540088871
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.522503666666667
$ws.Range("H5").Value = 4.567511000000001
$ws.Range("I5").Value = 0.2454047764062963
$ws.Range("J5").Value = 0.2454047764062963
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.044118333333333
$ws.Range("N5").Value = 6.132354999999999
$ws.Range("O5").Value = 0.1776005292722278
$ws.Range("P5").Value = 0.1776005292722278
$ws.Range("Q5").Value = 3.112177657600556
$ws.Range("R5").Value = 28.009598918405
$ws.Range("S5").Value = 0.04358401817569094
$ws.Range("T5").Value = 0.04358401817569093
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.522503666666667
$ws.Range("H6").Value = 4.567511000000001
$ws.Range("I6").Value = 0.2454047764062963
$ws.Range("J6").Value = 0.2454047764062963
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.059280333333334
$ws.Range("N6").Value = 21.177841
$ws.Range("O6").Value = 0.6133362746356149
$ws.Range("P6").Value = 0.6133362746356149
$ws.Range("Q6").Value = 10.74778019152789
$ws.Range("R6").Value = 96.73002172375101
$ws.Range("S6").Value = 0.1505156513388238
$ws.Range("T6").Value = 0.1505156513388238
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.522503666666667
$ws.Range("H7").Value = 4.567511000000001
$ws.Range("I7").Value = 0.2454047764062963
$ws.Range("J7").Value = 0.2454047764062963
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.406242333333333
$ws.Range("N7").Value = 7.218726999999999
$ws.Range("O7").Value = 0.2090631960921573
$ws.Range("P7").Value = 0.2090631960921573
$ws.Range("Q7").Value = 3.663512775388556
$ws.Range("R7").Value = 32.971614978497
$ws.Range("S7").Value = 0.05130510689178153
$ws.Range("T7").Value = 0.05130510689178153
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.7994406666666668
$ws.Range("H8").Value = 2.398322
$ws.Range("I8").Value = 0.1288578558782456
$ws.Range("J8").Value = 0.1288578558782456
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.044118333333333
$ws.Range("N8").Value = 6.132354999999999
$ws.Range("O8").Value = 0.1776005292722278
$ws.Range("P8").Value = 0.1776005292722278
$ws.Range("Q8").Value = 1.634151323145556
$ws.Range("R8").Value = 14.70736190831
$ws.Range("S8").Value = 0.02288522340486087
$ws.Range("T8").Value = 0.02288522340486086
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.7994406666666668
$ws.Range("H9").Value = 2.398322
$ws.Range("I9").Value = 0.1288578558782456
$ws.Range("J9").Value = 0.1288578558782456
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.059280333333334
$ws.Range("N9").Value = 21.177841
$ws.Range("O9").Value = 0.6133362746356149
$ws.Range("P9").Value = 0.6133362746356149
$ws.Range("Q9").Value = 5.64347577586689
$ws.Range("R9").Value = 50.791281982802
$ws.Range("S9").Value = 0.07903319728189613
$ws.Range("T9").Value = 0.07903319728189612
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.7994406666666668
$ws.Range("H10").Value = 2.398322
$ws.Range("I10").Value = 0.1288578558782456
$ws.Range("J10").Value = 0.1288578558782456
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.406242333333333
$ws.Range("N10").Value = 7.218726999999999
$ws.Range("O10").Value = 0.2090631960921573
$ws.Range("P10").Value = 0.2090631960921573
$ws.Range("Q10").Value = 1.923647975121556
$ws.Range("R10").Value = 17.312831776094
$ws.Range("S10").Value = 0.0269394351914886
$ws.Range("T10").Value = 0.0269394351914886
